$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly driver report update for 2025-04-21
$ws.Range("D3").Value = 94
$ws.Range("B4").Value = 55
$ws.Range("C4").Value = 1412
$ws.Range("B5").Value = 58
$ws.Range("C5").Value = 1637
$ws.Range("B13").Value = 449371
$ws.Range("B17").Value = 77999
